$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 454.41666
$ws.Cells.Item(80, 10).Value = 553.2222
$ws.Cells.Item(80, 12).Value = 1659.6666
$ws.Cells.Item(80, 14).Value = -3655.6666
$ws.Cells.Item(83, 8).Value = 454.41666
$ws.Cells.Item(83, 10).Value = 553.2222
$ws.Cells.Item(83, 12).Value = 4978.999800000001
$ws.Cells.Item(83, 14).Value = -14962.9998
$ws.Cells.Item(86, 8).Value = 4423
$ws.Cells.Item(86, 9).Value = 3948
$ws.Cells.Item(86, 10).Value = 5016.75
$ws.Cells.Item(86, 11).Value = 3948
$ws.Cells.Item(86, 12).Value = 5016.75
$ws.Cells.Item(86, 13).Value = -2825
$ws.Cells.Item(86, 14).Value = -7262.75
$ws.Cells.Item(89, 8).Value = 4423
$ws.Cells.Item(89, 9).Value = 3948
$ws.Cells.Item(89, 10).Value = 5016.75
$ws.Cells.Item(89, 11).Value = 19740
$ws.Cells.Item(89, 12).Value = 25083.75
$ws.Cells.Item(89, 13).Value = -14124
$ws.Cells.Item(89, 14).Value = -36315.75
$ws.Cells.Item(97, 8).Value = 28995.625
$ws.Cells.Item(97, 9).Value = 2983.3333
$ws.Cells.Item(97, 11).Value = 8949.999899999999
$ws.Cells.Item(97, 13).Value = -8453.999899999999
$ws.Cells.Item(100, 8).Value = 1732.0714
$ws.Cells.Item(100, 9).Value = 1995.091
$ws.Cells.Item(100, 10).Value = 767.6667
$ws.Cells.Item(100, 11).Value = 1995.091
$ws.Cells.Item(100, 12).Value = 767.6667
$ws.Cells.Item(100, 13).Value = -1454.091
$ws.Cells.Item(100, 14).Value = -1849.6667
$ws.Cells.Item(137, 8).Value = 1818.909
$ws.Cells.Item(137, 9).Value = 1835.7333
$ws.Cells.Item(137, 11).Value = 5507.199900000001
$ws.Cells.Item(137, 13).Value = -2957.199900000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6009.7466
$ws.Cells.Item(32, 9).Value = 3745.8462
$ws.Cells.Item(32, 10).Value = 30535.334
$ws.Cells.Item(32, 11).Value = 3745.8462
$ws.Cells.Item(32, 12).Value = 30535.334
$ws.Cells.Item(32, 13).Value = -3458.8462
$ws.Cells.Item(32, 14).Value = -31109.334
$ws.Cells.Item(61, 8).Value = 3830.5789
$ws.Cells.Item(61, 9).Value = 2196.6333
$ws.Cells.Item(61, 11).Value = 2196.6333
$ws.Cells.Item(61, 13).Value = -1984.6333
$ws.Cells.Item(69, 8).Value = 496940
$ws.Cells.Item(69, 10).Value = 496940
$ws.Cells.Item(69, 12).Value = 496940
$ws.Cells.Item(69, 14).Value = -498438
$ws.Cells.Item(72, 8).Value = 496940
$ws.Cells.Item(72, 10).Value = 496940
$ws.Cells.Item(72, 12).Value = 1490820
$ws.Cells.Item(72, 14).Value = -1498308
$ws.Cells.Item(97, 8).Value = 1532.0322
$ws.Cells.Item(97, 9).Value = 1358
$ws.Cells.Item(97, 10).Value = 1897.5
$ws.Cells.Item(97, 11).Value = 1358
$ws.Cells.Item(97, 12).Value = 1897.5
$ws.Cells.Item(97, 13).Value = -862
$ws.Cells.Item(97, 14).Value = -2889.5
$ws.Cells.Item(122, 8).Value = 1508.5416
$ws.Cells.Item(122, 10).Value = 1992.5555
$ws.Cells.Item(122, 12).Value = 5977.666499999999
$ws.Cells.Item(122, 14).Value = -10877.6665
$ws.Cells.Item(132, 8).Value = 2079.5454
$ws.Cells.Item(132, 9).Value = 1902.2709
$ws.Cells.Item(132, 11).Value = 5706.8127
$ws.Cells.Item(132, 13).Value = -3176.8127
$ws.Cells.Item(136, 8).Value = 3830.5789
$ws.Cells.Item(136, 9).Value = 2196.6333
$ws.Cells.Item(136, 11).Value = 6589.8999
$ws.Cells.Item(136, 13).Value = -4039.8999
$ws.Cells.Item(139, 8).Value = 107999.5
$ws.Cells.Item(139, 10).Value = 107999.5
$ws.Cells.Item(139, 12).Value = 107999.5
$ws.Cells.Item(139, 14).Value = -118279.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4500.6333
$ws.Cells.Item(20, 9).Value = 3506.6667
$ws.Cells.Item(20, 11).Value = 3506.6667
$ws.Cells.Item(20, 13).Value = -3259.6667
$ws.Cells.Item(86, 8).Value = 1266.1666
$ws.Cells.Item(86, 9).Value = 1239.3
$ws.Cells.Item(86, 10).Value = 1299.75
$ws.Cells.Item(86, 11).Value = 1239.3
$ws.Cells.Item(86, 12).Value = 1299.75
$ws.Cells.Item(86, 13).Value = -116.3
$ws.Cells.Item(86, 14).Value = -3545.75
$ws.Cells.Item(89, 8).Value = 1266.1666
$ws.Cells.Item(89, 9).Value = 1239.3
$ws.Cells.Item(89, 10).Value = 1299.75
$ws.Cells.Item(89, 11).Value = 6196.5
$ws.Cells.Item(89, 12).Value = 6498.75
$ws.Cells.Item(89, 13).Value = -580.5
$ws.Cells.Item(89, 14).Value = -17730.75
$ws.Cells.Item(107, 8).Value = 1280.8462
$ws.Cells.Item(107, 9).Value = 1259.5454
$ws.Cells.Item(107, 10).Value = 1398
$ws.Cells.Item(107, 11).Value = 1259.5454
$ws.Cells.Item(107, 12).Value = 1398
$ws.Cells.Item(107, 13).Value = 660.4546
$ws.Cells.Item(107, 14).Value = -5238
$ws.Cells.Item(134, 8).Value = 1527.9166
$ws.Cells.Item(134, 9).Value = 1491
$ws.Cells.Item(134, 11).Value = 4473
$ws.Cells.Item(134, 13).Value = -1938
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 39559.52
$ws.Cells.Item(31, 9).Value = 47629
$ws.Cells.Item(31, 11).Value = 47629
$ws.Cells.Item(31, 13).Value = -47334
$ws.Cells.Item(34, 8).Value = 39559.52
$ws.Cells.Item(34, 9).Value = 47629
$ws.Cells.Item(34, 11).Value = 47629
$ws.Cells.Item(34, 13).Value = -47427
$ws.Cells.Item(99, 8).Value = 3016.6667
$ws.Cells.Item(99, 9).Value = 2650
$ws.Cells.Item(99, 10).Value = 3750
$ws.Cells.Item(99, 11).Value = 2650
$ws.Cells.Item(99, 12).Value = 3750
$ws.Cells.Item(99, 13).Value = -1152
$ws.Cells.Item(99, 14).Value = -6746
$ws.Cells.Item(107, 8).Value = 1106.6666
$ws.Cells.Item(107, 9).Value = 875.3333
$ws.Cells.Item(107, 10).Value = 1222.3334
$ws.Cells.Item(107, 11).Value = 875.3333
$ws.Cells.Item(107, 12).Value = 1222.3334
$ws.Cells.Item(107, 13).Value = 1044.6667
$ws.Cells.Item(107, 14).Value = -5062.3334
$ws.Cells.Item(125, 8).Value = 57326
$ws.Cells.Item(125, 10).Value = 57326
$ws.Cells.Item(125, 12).Value = 57326
$ws.Cells.Item(125, 14).Value = -62246
$ws.Cells.Item(126, 8).Value = 3016.6667
$ws.Cells.Item(126, 9).Value = 2650
$ws.Cells.Item(126, 10).Value = 3750
$ws.Cells.Item(126, 11).Value = 7950
$ws.Cells.Item(126, 12).Value = 11250
$ws.Cells.Item(126, 13).Value = -5480
$ws.Cells.Item(126, 14).Value = -16190
$ws.Cells.Item(134, 8).Value = 15595.883
$ws.Cells.Item(134, 9).Value = 5169.7095
$ws.Cells.Item(134, 10).Value = 123333
$ws.Cells.Item(134, 11).Value = 15509.1285
$ws.Cells.Item(134, 12).Value = 369999
$ws.Cells.Item(134, 13).Value = -12974.1285
$ws.Cells.Item(134, 14).Value = -375069
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8058
$ws.Cells.Item(70, 9).Value = 6869.6
$ws.Cells.Item(70, 10).Value = 14000
$ws.Cells.Item(70, 11).Value = 6869.6
$ws.Cells.Item(70, 12).Value = 14000
$ws.Cells.Item(70, 13).Value = -6599.6
$ws.Cells.Item(70, 14).Value = -14540
$ws.Cells.Item(73, 8).Value = 8058
$ws.Cells.Item(73, 9).Value = 6869.6
$ws.Cells.Item(73, 10).Value = 14000
$ws.Cells.Item(73, 11).Value = 6869.6
$ws.Cells.Item(73, 12).Value = 14000
$ws.Cells.Item(73, 13).Value = -5933.6
$ws.Cells.Item(73, 14).Value = -15872
$ws.Cells.Item(102, 8).Value = 90910170
$ws.Cells.Item(102, 9).Value = 867.2222
$ws.Cells.Item(102, 11).Value = 867.2222
$ws.Cells.Item(102, 13).Value = 754.7778
$ws.Cells.Item(122, 8).Value = 3235.4443
$ws.Cells.Item(122, 9).Value = 2794.3667
$ws.Cells.Item(122, 10).Value = 5440.8335
$ws.Cells.Item(122, 11).Value = 8383.1001
$ws.Cells.Item(122, 12).Value = 16322.5005
$ws.Cells.Item(122, 13).Value = -5933.1001
$ws.Cells.Item(122, 14).Value = -21222.5005
$ws.Cells.Item(132, 8).Value = 3042.7021
$ws.Cells.Item(132, 9).Value = 2533.3333
$ws.Cells.Item(132, 10).Value = 4709.727
$ws.Cells.Item(132, 11).Value = 7599.999899999999
$ws.Cells.Item(132, 12).Value = 14129.181
$ws.Cells.Item(132, 13).Value = -5069.999899999999
$ws.Cells.Item(132, 14).Value = -19189.181
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2562.5652
$ws.Cells.Item(16, 9).Value = 2664.1428
$ws.Cells.Item(16, 11).Value = 2664.1428
$ws.Cells.Item(16, 13).Value = -2494.1428
$ws.Cells.Item(40, 8).Value = 3218.4
$ws.Cells.Item(40, 9).Value = 2425.1155
$ws.Cells.Item(40, 11).Value = 2425.1155
$ws.Cells.Item(40, 13).Value = -2289.1155
$ws.Cells.Item(46, 8).Value = 2671.125
$ws.Cells.Item(46, 9).Value = 1873.8
$ws.Cells.Item(46, 11).Value = 1873.8
$ws.Cells.Item(46, 13).Value = -1685.8
$ws.Cells.Item(61, 8).Value = 4685.55
$ws.Cells.Item(61, 9).Value = 4537.533
$ws.Cells.Item(61, 11).Value = 4537.533
$ws.Cells.Item(61, 13).Value = -4335.533
$ws.Cells.Item(113, 8).Value = 4685.55
$ws.Cells.Item(113, 9).Value = 4537.533
$ws.Cells.Item(113, 11).Value = 4537.533
$ws.Cells.Item(113, 13).Value = -2367.533
$ws.Cells.Item(136, 8).Value = 5381.591
$ws.Cells.Item(136, 9).Value = 4523.2354
$ws.Cells.Item(136, 10).Value = 8300
$ws.Cells.Item(136, 11).Value = 13569.7062
$ws.Cells.Item(136, 12).Value = 24900
$ws.Cells.Item(136, 13).Value = -11019.7062
$ws.Cells.Item(136, 14).Value = -30000
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 14249.25
$ws.Cells.Item(15, 10).Value = 18498.5
$ws.Cells.Item(15, 12).Value = 18498.5
$ws.Cells.Item(15, 14).Value = -19074.5
$ws.Cells.Item(126, 8).Value = 2142.5
$ws.Cells.Item(126, 9).Value = 2049.7693
$ws.Cells.Item(126, 10).Value = 2314.7144
$ws.Cells.Item(126, 11).Value = 6149.3079
$ws.Cells.Item(126, 12).Value = 6944.1432
$ws.Cells.Item(126, 13).Value = -3679.3079
$ws.Cells.Item(126, 14).Value = -11884.1432
$ws.Cells.Item(132, 8).Value = 1802.5938
$ws.Cells.Item(132, 9).Value = 1722.7667
$ws.Cells.Item(132, 11).Value = 5168.300099999999
$ws.Cells.Item(132, 13).Value = -2638.300099999999

Write-Output "applied 222 changes"